$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update StatQuery text (shared by C2:C5) - set on C2 first so it becomes the new shared string in place of the old one
$ws.Range("C2").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (f:file)-[*]->(samp:sample)-->(c)
MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp,demo, c, s, p, diag
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Airedale Terrier'', ''Bluetick Hound'', ''Norfolk Terrier'', ''Norwegian Elkhound'', ''Scottish Terrier'', ''Welsh Springer Spaniel'', ''Wheaten Terrier'']and diag.disease_term in [''Bladder Cancer''] and demo.sex in [''Male'', ''Female''] and demo.neutered_indicator IN [''No'', ''Yes'']
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

# 2) New row 5 tab-name cell introduces the new tab label next
$ws.Range("A5").Value = 'StudyFilesTab'

# 3) Reworked CasesTab query (B2)
$ws.Range("B2").Value = ' MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Airedale Terrier'', ''Bluetick Hound'', ''Norfolk Terrier'', ''Norwegian Elkhound'', ''Scottish Terrier'', ''Welsh Springer Spaniel'', ''Wheaten Terrier'']and diag.disease_term in [''Bladder Cancer''] and demo.sex in [''Male'', ''Female''] and demo.neutered_indicator IN [''No'', ''Yes'']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '''') AS `Case ID`,
       coalesce(s.clinical_study_designation, '''') AS `Study Code`,
       coalesce(s.clinical_study_type, '''') AS  `Study Type`,
       coalesce(demo.breed, '''') AS Breed ,
       coalesce(diag.disease_term, '''') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '''') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '''') AS Sex,
       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,
       coalesce(diag.best_response, '''') AS `Response to Treatment`,
       coalesce(co.cohort_description, '''') AS `Cohort`
Order by c.case_id LIMIT 100        '

# 4) Reworked SamplesTab query (B3)
$ws.Range("B3").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Airedale Terrier'', ''Bluetick Hound'', ''Norfolk Terrier'', ''Norwegian Elkhound'', ''Scottish Terrier'', ''Welsh Springer Spaniel'', ''Wheaten Terrier'']and diag.disease_term in [''Bladder Cancer''] and demo.sex in [''Male'', ''Female''] and demo.neutered_indicator IN [''No'', ''Yes'']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'

# 5) Reworked FilesTab query (B4)
$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Airedale Terrier'', ''Bluetick Hound'', ''Norfolk Terrier'', ''Norwegian Elkhound'', ''Scottish Terrier'', ''Welsh Springer Spaniel'', ''Wheaten Terrier'']and diag.disease_term in [''Bladder Cancer''] and demo.sex in [''Male'', ''Female''] and demo.neutered_indicator IN [''No'', ''Yes'']
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '''') AS `File Name`, 
 coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_type, '''') AS `File Type`, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
   coalesce(samp.sample_id, '''') AS `Sample ID`,
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis 
        Order By f.file_name LIMIT 100'

# 6) New StudyFilesTab query (B5)
$ws.Range("B5").Value = 'MATCH (f:file)-->(s:study)
MATCH (s)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (sf:file)-->(s)
MATCH (s)<--(c)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN [''UBC02''] and demo.breed in [''Airedale Terrier'', ''Bluetick Hound'', ''Norfolk Terrier'', ''Norwegian Elkhound'', ''Scottish Terrier'', ''Welsh Springer Spaniel'', ''Wheaten Terrier'']and diag.disease_term in [''Bladder Cancer''] and demo.sex in [''Male'', ''Female''] and demo.neutered_indicator IN [''No'', ''Yes'']
WITH DISTINCT f,  s, c
WITH
        f, c,  s,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c,  s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c,   s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '''') AS `File Name`,
  coalesce(f.file_type, '''') AS `File Type`,
  coalesce("study", '''') AS `Association`,
  coalesce(f.file_description, '''') AS `Description`,
  coalesce(f.file_format, '''') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
  coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# 7) Fill the rest of the StatQuery column for the other tabs + new row
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("C5").Value = $ws.Range("C2").Value2

# 8) Fill D5/E5 with the same workbook file-name references as the other rows
$ws.Range("D5").Value = $ws.Range("D4").Value2
$ws.Range("E5").Value = $ws.Range("E4").Value2

# 9) Match wrap formatting used by the other query/stat columns
$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true

# 10) Row heights to fit the new/edited text
$ws.Rows.Item(2).RowHeight = 390
$ws.Rows.Item(3).RowHeight = 315
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# 11) Update the active selection to the new row
$ws.Range("C5").Select()

